# Add a new risk register row (row 5) to Sheet1 and update the scatter chart
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Value = "I have to present my demo early"
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 20
$ws.Range("D5").Formula = "=B5*C5"
$ws.Range("E5").Value = "4/20/2023"
$ws.Range("G5").Value = "O"
$ws.Range("H5").Value = "I will plan to get ahead in case my demo is early"

$ws.Range("H14").Select()
